$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# IDs in column A cycle through these 6 values (same pattern as existing rows)
$ids = @("i076652", "c5169789", "d2122378", "i312231", "i066088", "i047892")

# New names for the "BD" batch (rows 25-30) and "CE" batch (rows 31-36)
$namesBD = @("Zhang, SanBD", "Li, SiBd", "Wang, WuBD", "Zhou, LiuBD", "Wu, QiBD", "Zheng, BaBD")
$namesCE = @("Zhang, SanCE", "Li, SiCE", "Wang, WuCE", "Zhou, LiuCE", "Wu, QiCE", "Zheng, BaCE")

$row = 25
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $namesBD[$i]
    $row++
}
for ($i = 0; $i -lt 6; $i++) {
    $ws.Cells.Item($row, 1).Value = $ids[$i]
    $ws.Cells.Item($row, 2).Value = $namesCE[$i]
    $row++
}

$ws.Columns.Item(2).ColumnWidth = 24.79

$ws.Range("B35").Select()
